# Weekly update to the Ají (Hortaliza) price sheet:
#   - two new price records (row 233 & 234) are inserted at the top of the
#     "Terminal La Palmera de La Serena" Ají block, pushing the existing
#     records (old rows 233-242) down by two rows (to 235-244).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 233 - this shifts the existing data
# (old rows 233-242) down to rows 235-244 and carries formatting/styles
# down with it (matches Excel's native "Insert Row" behaviour).
$ws.Rows.Item(233).Insert()
$ws.Rows.Item(233).Insert()

# --- New row 233 -----------------------------------------------------
$ws.Cells.Item(233, 1).Value2 = 8
$ws.Cells.Item(233, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(233, 3).Value = "Coquimbo"
$ws.Cells.Item(233, 4).Value2 = 44753
$ws.Cells.Item(233, 5).Value2 = 4
$ws.Cells.Item(233, 6).Value2 = 100112021
$ws.Cells.Item(233, 7).Value = "Ají"
$ws.Cells.Item(233, 8).Value = "Inferno"
$ws.Cells.Item(233, 9).Value = "Primera"
$ws.Cells.Item(233, 10).Value2 = 480
$ws.Cells.Item(233, 11).Value2 = 14000
$ws.Cells.Item(233, 12).Value2 = 15000
$ws.Cells.Item(233, 13).Value2 = 14500
$ws.Cells.Item(233, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(233, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(233, 16).Value2 = 1208
$ws.Cells.Item(233, 17).Value2 = 12
$ws.Cells.Item(233, 18).Value = "Hortaliza"

# --- New row 234 -----------------------------------------------------
$ws.Cells.Item(234, 1).Value2 = 8
$ws.Cells.Item(234, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(234, 3).Value = "Coquimbo"
$ws.Cells.Item(234, 4).Value2 = 44753
$ws.Cells.Item(234, 5).Value2 = 4
$ws.Cells.Item(234, 6).Value2 = 100112021
$ws.Cells.Item(234, 7).Value = "Ají"
$ws.Cells.Item(234, 8).Value = "Inferno"
$ws.Cells.Item(234, 9).Value = "Segunda"
$ws.Cells.Item(234, 10).Value2 = 200
$ws.Cells.Item(234, 11).Value2 = 9000
$ws.Cells.Item(234, 12).Value2 = 10000
$ws.Cells.Item(234, 13).Value2 = 9500
$ws.Cells.Item(234, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(234, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(234, 16).Value2 = 792
$ws.Cells.Item(234, 17).Value2 = 12
$ws.Cells.Item(234, 18).Value = "Hortaliza"
